$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: write new text labels in the exact order they must be
# appended to sharedStrings.xml (matches the target uniqueCount ordering) ---
$ws.Range("A13").Value = "PLL"
$ws.Range("A14").Value = "SCK"
$ws.Range("A16").Value = "K"
$ws.Range("A19").Value = "PLLCK"
$ws.Range("A20").Value = "Fs ref"
$ws.Range("A22").Value = "0x25 div"
$ws.Range("A23").Value = "0x26 div"
$ws.Range("A15").Value = "0x2A R"
$ws.Range("A17").Value = "0x29 P"
$ws.Range("A24").Value = "BCK"
$ws.Range("A25").Value = "0x27 div"
$ws.Range("A26").Value = "LRCK"
$ws.Range("A27").Value = "BCK/Fs"

# --- Phase 2: remaining values (numbers / formulas / reused labels) ---
$ws.Range("A18").Value = "N"
$ws.Range("B18").Value = 2048

$ws.Range("B14").Value = 8000
$ws.Range("C14").Value = "kHz"

$ws.Range("B15").Value = 1

$ws.Range("B16").Value = 12

$ws.Range("B17").Value = 1

$ws.Range("B19").Formula = "=B14*B15*B16/B17"
$ws.Range("C19").Value = "kHz"

$ws.Range("B20").Formula = "=B19/B18"
$ws.Range("C20").Value = "kHz"

$ws.Range("B22").Value = 8

$ws.Range("B23").Value = 4

$ws.Range("B24").Formula = "=B19/(B22*B23)"
$ws.Range("C24").Value = "kHz"

$ws.Range("B25").Value = 64

$ws.Range("B26").Formula = "=B19/(B22*B23*B25)"

$ws.Range("B27").Formula = "=B24/B20"

# --- Phase 3: styles ---
$ws.Range("A13").Style = "Accent1"

$ws.Range("A14:B14").Style = "Good"
$ws.Range("A15:B15").Style = "Good"
$ws.Range("A16:B16").Style = "Good"
$ws.Range("A17:B17").Style = "Good"

$ws.Range("A19:B19").Style = "Calculation"
$ws.Range("A20:B20").Style = "Calculation"

$ws.Range("A22:B22").Style = "Good"
$ws.Range("A23:B23").Style = "Good"

$ws.Range("A24:B24").Style = "Calculation"

$ws.Range("A25:B25").Style = "Good"

$ws.Range("A26:B26").Style = "Calculation"
$ws.Range("A27:B27").Style = "Calculation"

# --- Phase 4: selection, matching the saved view state ---
$ws.Range("B23").Select()
